$d = $word.ActiveDocument

$replacements = @(
    @("92×39=", "18×18="),
    @("98×79=", "95×88="),
    @("74×60=", "14×83="),
    @("83×32=", "65×96="),
    @("81×53=", "87×60="),
    @("20×92=", "71×56="),
    @("65×54=", "15×57="),
    @("97×35=", "76×34="),
    @("60×43=", "72×63="),
    @("57×53=", "60×11="),
    @("29×15=", "68×74="),
    @("38×46=", "63×19="),
    @("91×63=", "66×29="),
    @("33×61=", "49×74="),
    @("38×67=", "95×56="),
    @("47×88=", "42×49="),
    @("65×22=", "71×22="),
    @("55×35=", "91×96="),
    @("56×23=", "98×78="),
    @("52×60=", "26×59="),
    @("43×17=", "78×43="),
    @("33×51=", "58×70="),
    @("53×31=", "51×97="),
    @("77×38=", "37×34="),
    @("18×68=", "37×55=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
